{"js": "// Replace the header date and the 25 multiplication problems with their\n// updated values. Every old string in this worksheet is unique, so a\n// plain case-sensitive search-and-replace on each pair is safe and keeps\n// all run formatting (fonts/size) untouched since we only touch the text\n// of the matched range.\nconst replacements = [\n  [\"2026-01-07 Wednesday\", \"2026-01-08 Thursday\"],\n  [\"480\u00d77=\", \"427\u00d79=\"],\n  [\"221\u00d78=\", \"882\u00d74=\"],\n  [\"590\u00d75=\", \"835\u00d78=\"],\n  [\"706\u00d73=\", \"637\u00d78=\"],\n  [\"459\u00d78=\", \"491\u00d79=\"],\n  [\"395\u00d76=\", \"620\u00d77=\"],\n  [\"726\u00d72=\", \"142\u00d76=\"],\n  [\"924\u00d79=\", \"964\u00d72=\"],\n  [\"340\u00d73=\", \"528\u00d75=\"],\n  [\"817\u00d77=\", \"677\u00d72=\"],\n  [\"332\u00d75=\", \"847\u00d73=\"],\n  [\"869\u00d79=\", \"939\u00d73=\"],\n  [\"399\u00d74=\", \"812\u00d72=\"],\n  [\"836\u00d74=\", \"286\u00d78=\"],\n  [\"612\u00d76=\", \"283\u00d76=\"],\n  [\"629\u00d79=\", \"690\u00d76=\"],\n  [\"152\u00d74=\", \"996\u00d74=\"],\n  [\"951\u00d75=\", \"911\u00d76=\"],\n  [\"879\u00d77=\", \"771\u00d78=\"],\n  [\"646\u00d75=\", \"395\u00d79=\"],\n  [\"338\u00d74=\", \"983\u00d79=\"],\n  [\"321\u00d76=\", \"240\u00d79=\"],\n  [\"562\u00d75=\", \"526\u00d74=\"],\n  [\"847\u00d72=\", \"700\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the header date and the 25 multiplication problems with their\n# updated values. Every \"old\" string in this worksheet occurs exactly\n# once in the document, so a straightforward Find/Execute replace (no\n# wildcards) for each pair is safe and leaves all other formatting\n# (fonts/size/alignment) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Find = \"2026-01-07 Wednesday\"; Replace = \"2026-01-08 Thursday\" },\n    @{ Find = \"480\u00d77=\"; Replace = \"427\u00d79=\" },\n    @{ Find = \"221\u00d78=\"; Replace = \"882\u00d74=\" },\n    @{ Find = \"590\u00d75=\"; Replace = \"835\u00d78=\" },\n    @{ Find = \"706\u00d73=\"; Replace = \"637\u00d78=\" },\n    @{ Find = \"459\u00d78=\"; Replace = \"491\u00d79=\" },\n    @{ Find = \"395\u00d76=\"; Replace = \"620\u00d77=\" },\n    @{ Find = \"726\u00d72=\"; Replace = \"142\u00d76=\" },\n    @{ Find = \"924\u00d79=\"; Replace = \"964\u00d72=\" },\n    @{ Find = \"340\u00d73=\"; Replace = \"528\u00d75=\" },\n    @{ Find = \"817\u00d77=\"; Replace = \"677\u00d72=\" },\n    @{ Find = \"332\u00d75=\"; Replace = \"847\u00d73=\" },\n    @{ Find = \"869\u00d79=\"; Replace = \"939\u00d73=\" },\n    @{ Find = \"399\u00d74=\"; Replace = \"812\u00d72=\" },\n    @{ Find = \"836\u00d74=\"; Replace = \"286\u00d78=\" },\n    @{ Find = \"612\u00d76=\"; Replace = \"283\u00d76=\" },\n    @{ Find = \"629\u00d79=\"; Replace = \"690\u00d76=\" },\n    @{ Find = \"152\u00d74=\"; Replace = \"996\u00d74=\" },\n    @{ Find = \"951\u00d75=\"; Replace = \"911\u00d76=\" },\n    @{ Find = \"879\u00d77=\"; Replace = \"771\u00d78=\" },\n    @{ Find = \"646\u00d75=\"; Replace = \"395\u00d79=\" },\n    @{ Find = \"338\u00d74=\"; Replace = \"983\u00d79=\" },\n    @{ Find = \"321\u00d76=\"; Replace = \"240\u00d79=\" },\n    @{ Find = \"562\u00d75=\"; Replace = \"526\u00d74=\" },\n    @{ Find = \"847\u00d72=\"; Replace = \"700\u00d73=\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n}\n"}
